# Populate the "AI", "IoT", and "ハプティクス" sheets with scraped
# Wired.jp article URLs, and make the "AI" sheet the active/selected sheet
# (matching the state the workbook was left in after the scrape finished).

$wb = $excel.ActiveWorkbook

$aiUrls = @(
    "https://wired.jp/2019/12/06/qualcomm-snapdragon-865-5g/",
    "https://wired.jp/waia/2019/17_kiyoko-ojima/",
    "https://wired.jp/membership/2019/12/05/is-amazon-unstoppable/",
    "https://wired.jp/membership/2019/12/03/how-to-practice-long-term-thinking/",
    "https://wired.jp/2019/12/03/apple-china-iphone-sales-hong-kong-protests/",
    "https://wired.jp/series/away-from-animals-and-machines/",
    "https://wired.jp/series/away-from-animals-and-machines/chapter12-4/",
    "https://wired.jp/membership/2019/12/02/hidden-costs-automated-thinking/",
    "https://wired.jp/2019/11/30/los-angeles-blade-runner-theory/",
    "https://wired.jp/waia/2019/16_yukiko-ogawa/",
    "https://wired.jp/membership/2019/11/28/will-artificial-intelligence3/",
    "https://wired.jp/2019/11/27/viral-app-labels-you-isnt-what-you-think/",
    "https://wired.jp/membership/2019/11/21/will-artificial-intelligence2/",
    "https://wired.jp/2019/11/15/ars-hakuhodo-ws/",
    "https://wired.jp/2019/11/15/microsoft-sends-a-new-kind-of-ai-processor-into-the-cloud/"
)

$iotUrls = @(
    "https://wired.jp/event/kyotouniversity-vol35/",
    "https://wired.jp/2019/12/06/qualcomm-snapdragon-865-5g/",
    "https://wired.jp/2019/12/05/amazon-joins-quantum-computing-race/",
    "https://wired.jp/2019/12/05/larry-page-sergey-brin-step-down/",
    "https://wired.jp/culture/",
    "https://wired.jp/news/",
    "https://wired.jp/2019/12/03/lil-bub-cat-obit/",
    "https://wired.jp/2019/12/03/airmega-ws/",
    "https://wired.jp/2019/12/02/brewdog-hybrid-plant-beef-burger/",
    "https://wired.jp/2019/12/02/tiktok-time/",
    "https://wired.jp/2019/12/02/why-lightning-strikes-twice-as-much-over-shipping-lanes/"
)

$hapticsUrls = @(
    "https://wired.jp/2019/11/05/oracle-ana-ws/",
    "https://wired.jp/2019/10/09/exclusive-playstation-5/",
    "https://wired.jp/2019/06/14/how-i-became-a-robot-in-london/",
    "https://wired.jp/2018/10/06/native-instruments-traktor-pro-3/",
    "https://wired.jp/2018/04/21/ready-player-one-movie-vs-book/",
    "https://wired.jp/2016/10/09/ceatec-japan-2016-sensing/",
    "https://wired.jp/innovationinsights/post/wired/w/new-economy/",
    "https://wired.jp/2015/04/12/apples-haptic-tech/",
    "https://wired.jp/2015/01/19/next-world-06/",
    "https://wired.jp/2003/07/04/%E3%82%A4%E3%83%B3%E3%82%BF%E3%83%BC%E3%83%8D%E3%83%83%E3%83%88%E4%B8%8A%E3%81%A7%E8%A7%A6%E8%A6%9A%E3%82%92%E5%85%B1%E6%9C%89%E3%81%A7%E3%81%8D%E3%82%8B%E3%80%8E%E3%83%8F%E3%83%97%E3%83%86%E3%82%A3/",
    "https://wired.jp/2001/08/22/%E8%A7%A6%E8%A6%9A%E3%82%A4%E3%83%B3%E3%82%BF%E3%83%BC%E3%83%95%E3%82%A7%E3%83%BC%E3%82%B9%E3%80%8E%E3%83%8F%E3%83%97%E3%83%86%E3%82%A3%E3%83%83%E3%82%AF%E3%82%B9%E3%80%8F%E3%81%AE%E5%8F%AF%E8%83%BD/"
)

$wsAI = $wb.Worksheets.Item(1)
for ($i = 0; $i -lt $aiUrls.Length; $i++) {
    $wsAI.Cells.Item($i + 3, 1).Value = $aiUrls[$i]
}

$wsIoT = $wb.Worksheets.Item(2)
for ($i = 0; $i -lt $iotUrls.Length; $i++) {
    $wsIoT.Cells.Item($i + 3, 1).Value = $iotUrls[$i]
}

$wsHaptics = $wb.Worksheets.Item(3)
for ($i = 0; $i -lt $hapticsUrls.Length; $i++) {
    $wsHaptics.Cells.Item($i + 3, 1).Value = $hapticsUrls[$i]
}

# Move the active tab / selection from the last sheet ("VR広告") to the
# first sheet ("AI"), landing the cursor just past the newly written data.
$wsAI.Activate() | Out-Null
$wsAI.Range("A25").Select() | Out-Null
